$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 165.66667
$ws.Range("I12").Value = 98
$ws.Range("J12").Value = 199.5
$ws.Range("K12").Value = 98
$ws.Range("L12").Value = 199.5
$ws.Range("M12").Value = 72
$ws.Range("N12").Value = -539.5
$ws.Range("H41").Value = 1571.4286
$ws.Range("I41").Value = 1343.8
$ws.Range("K41").Value = 1343.8
$ws.Range("M41").Value = -903.8
$ws.Range("H80").Value = 722.9545000000001
$ws.Range("I80").Value = 584.06665
$ws.Range("J80").Value = 1020.5714
$ws.Range("K80").Value = 1752.19995
$ws.Range("L80").Value = 3061.7142
$ws.Range("M80").Value = -754.1999499999999
$ws.Range("N80").Value = -5057.7142
$ws.Range("H83").Value = 722.9545000000001
$ws.Range("I83").Value = 584.06665
$ws.Range("J83").Value = 1020.5714
$ws.Range("K83").Value = 5256.59985
$ws.Range("L83").Value = 9185.142600000001
$ws.Range("M83").Value = -264.5998499999996
$ws.Range("N83").Value = -19169.1426
$ws.Range("H112").Value = 1911865.8
$ws.Range("J112").Value = 2210465.5
$ws.Range("L112").Value = 6631396.5
$ws.Range("N112").Value = -6633612.5
$ws.Range("H137").Value = 6068.0884
$ws.Range("I137").Value = 7724.6875
$ws.Range("J137").Value = 4595.5557
$ws.Range("K137").Value = 23174.0625
$ws.Range("L137").Value = 13786.6671
$ws.Range("M137").Value = -20624.0625
$ws.Range("N137").Value = -18886.6671
$ws.Range("H138").Value = 2774.15
$ws.Range("I138").Value = 1875.875
$ws.Range("J138").Value = 3373
$ws.Range("K138").Value = 5627.625
$ws.Range("L138").Value = 10119
$ws.Range("M138").Value = -487.625
$ws.Range("N138").Value = -20399
$ws.Range("H141").Value = 4662.5557
$ws.Range("I141").Value = 5127.3335
$ws.Range("K141").Value = 15382.0005
$ws.Range("M141").Value = -10202.0005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26997.072
$ws.Range("I32").Value = 26126.88
$ws.Range("J32").Value = 34248.668
$ws.Range("K32").Value = 26126.88
$ws.Range("L32").Value = 34248.668
$ws.Range("M32").Value = -25839.88
$ws.Range("N32").Value = -34822.668
$ws.Range("H61").Value = 5105.3887
$ws.Range("I61").Value = 2136
$ws.Range("J61").Value = 8817.125
$ws.Range("K61").Value = 2136
$ws.Range("L61").Value = 8817.125
$ws.Range("M61").Value = -1924
$ws.Range("N61").Value = -9241.125
$ws.Range("H62").Value = 59800
$ws.Range("J62").Value = 59800
$ws.Range("L62").Value = 59800
$ws.Range("N62").Value = -61048
$ws.Range("H65").Value = 59800
$ws.Range("J65").Value = 59800
$ws.Range("L65").Value = 179400
$ws.Range("N65").Value = -185640
$ws.Range("H74").Value = 3929.1018
$ws.Range("I74").Value = 3982.362
$ws.Range("J74").Value = 840
$ws.Range("K74").Value = 3982.362
$ws.Range("L74").Value = 840
$ws.Range("M74").Value = -3108.362
$ws.Range("N74").Value = -2588
$ws.Range("H77").Value = 3929.1018
$ws.Range("I77").Value = 3982.362
$ws.Range("J77").Value = 840
$ws.Range("K77").Value = 19911.81
$ws.Range("L77").Value = 4200
$ws.Range("M77").Value = -15543.81
$ws.Range("N77").Value = -12936
$ws.Range("H110").Value = 2818.5
$ws.Range("I110").Value = 1575.909
$ws.Range("J110").Value = 6235.625
$ws.Range("K110").Value = 1575.909
$ws.Range("L110").Value = 6235.625
$ws.Range("M110").Value = 469.0909999999999
$ws.Range("N110").Value = -10325.625
$ws.Range("H122").Value = 5655.927
$ws.Range("J122").Value = 30278.666
$ws.Range("L122").Value = 90835.99800000001
$ws.Range("N122").Value = -95735.99800000001
$ws.Range("H136").Value = 5105.3887
$ws.Range("I136").Value = 2136
$ws.Range("J136").Value = 8817.125
$ws.Range("K136").Value = 6408
$ws.Range("L136").Value = 26451.375
$ws.Range("M136").Value = -3858
$ws.Range("N136").Value = -31551.375

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 54000
$ws.Range("J63").Value = 54000
$ws.Range("L63").Value = 54000
$ws.Range("N63").Value = -55372
$ws.Range("H66").Value = 54000
$ws.Range("J66").Value = 54000
$ws.Range("L66").Value = 162000
$ws.Range("N66").Value = -168864
$ws.Range("H96").Value = 3000
$ws.Range("I96").Value = 3000
$ws.Range("K96").Value = 3000
$ws.Range("M96").Value = -254
$ws.Range("H105").Value = 1524.7931
$ws.Range("I105").Value = 1457.1666
$ws.Range("J105").Value = 1849.4
$ws.Range("K105").Value = 1457.1666
$ws.Range("L105").Value = 1849.4
$ws.Range("M105").Value = 289.8334
$ws.Range("N105").Value = -5343.4
$ws.Range("H107").Value = 2821.0881
$ws.Range("I107").Value = 2726.1072
$ws.Range("J107").Value = 3264.3333
$ws.Range("K107").Value = 2726.1072
$ws.Range("L107").Value = 3264.3333
$ws.Range("M107").Value = -806.1071999999999
$ws.Range("N107").Value = -7104.3333
$ws.Range("H134").Value = 2011
$ws.Range("I134").Value = 1859.1333
$ws.Range("K134").Value = 5577.3999
$ws.Range("M134").Value = -3042.3999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18884864
$ws.Range("I31").Value = 1417.3636
$ws.Range("J31").Value = 50042548
$ws.Range("K31").Value = 1417.3636
$ws.Range("L31").Value = 50042548
$ws.Range("M31").Value = -1122.3636
$ws.Range("N31").Value = -50043138
$ws.Range("H34").Value = 18884864
$ws.Range("I34").Value = 1417.3636
$ws.Range("J34").Value = 50042548
$ws.Range("K34").Value = 1417.3636
$ws.Range("L34").Value = 50042548
$ws.Range("M34").Value = -1215.3636
$ws.Range("N34").Value = -50042952
$ws.Range("H86").Value = 33338498
$ws.Range("J86").Value = 6530.625
$ws.Range("L86").Value = 6530.625
$ws.Range("N86").Value = -8776.625
$ws.Range("H89").Value = 33338498
$ws.Range("J89").Value = 6530.625
$ws.Range("L89").Value = 32653.125
$ws.Range("N89").Value = -43885.125
$ws.Range("H107").Value = 4099.114
$ws.Range("I107").Value = 449.22223
$ws.Range("J107").Value = 6625.9614
$ws.Range("K107").Value = 449.22223
$ws.Range("L107").Value = 6625.9614
$ws.Range("M107").Value = 1470.77777
$ws.Range("N107").Value = -10465.9614

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 14308
$ws.Range("J32").Value = 14308
$ws.Range("L32").Value = 42924
$ws.Range("N32").Value = -43490
$ws.Range("H46").Value = 6875957
$ws.Range("J46").Value = 7858173.5
$ws.Range("L46").Value = 23574520.5
$ws.Range("N46").Value = -23574702.5
$ws.Range("H56").Value = 9792.214
$ws.Range("I56").Value = 9792.214
$ws.Range("K56").Value = 9792.214
$ws.Range("M56").Value = -9262.214
$ws.Range("H86").Value = 607.6
$ws.Range("I86").Value = 607.6
$ws.Range("K86").Value = 1822.8
$ws.Range("M86").Value = -636.8000000000002
$ws.Range("H89").Value = 607.6
$ws.Range("I89").Value = 607.6
$ws.Range("K89").Value = 5468.400000000001
$ws.Range("M89").Value = 459.5999999999995
$ws.Range("H140").Value = 5954054
$ws.Range("I140").Value = 19231770
$ws.Range("J140").Value = 1974.3448
$ws.Range("K140").Value = 57695310
$ws.Range("L140").Value = 5923.0344
$ws.Range("M140").Value = -57690130
$ws.Range("N140").Value = -16283.0344

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1432.5
$ws.Range("I31").Value = 1432.5
$ws.Range("K31").Value = 1432.5
$ws.Range("M31").Value = -1140.5
$ws.Range("H37").Value = 1432.5
$ws.Range("I37").Value = 1432.5
$ws.Range("K37").Value = 1432.5
$ws.Range("M37").Value = -1155.5
$ws.Range("H70").Value = 71438310
$ws.Range("I70").Value = 8049.75
$ws.Range("J70").Value = 166678670
$ws.Range("K70").Value = 8049.75
$ws.Range("L70").Value = 166678670
$ws.Range("M70").Value = -7779.75
$ws.Range("N70").Value = -166679210
$ws.Range("H73").Value = 71438310
$ws.Range("I73").Value = 8049.75
$ws.Range("J73").Value = 166678670
$ws.Range("K73").Value = 8049.75
$ws.Range("L73").Value = 166678670
$ws.Range("M73").Value = -7113.75
$ws.Range("N73").Value = -166680542
$ws.Range("H122").Value = 1386.7826
$ws.Range("I122").Value = 1185.5714
$ws.Range("K122").Value = 3556.7142
$ws.Range("M122").Value = -1106.7142

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1789.762
$ws.Range("I16").Value = 918.04877
$ws.Range("K16").Value = 918.04877
$ws.Range("M16").Value = -748.04877
$ws.Range("H122").Value = 3593.4666
$ws.Range("I122").Value = 3274.9
$ws.Range("K122").Value = 9824.700000000001
$ws.Range("M122").Value = -7374.700000000001
$ws.Range("H123").Value = 29999.5
$ws.Range("H136").Value = 2237.7114
$ws.Range("I136").Value = 1825.7273
$ws.Range("J136").Value = 4503.625
$ws.Range("K136").Value = 5477.1819
$ws.Range("L136").Value = 13510.875
$ws.Range("M136").Value = -2927.1819
$ws.Range("N136").Value = -18610.875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 16665.666
$ws.Range("H52").Value = 37132.832
$ws.Range("J52").Value = 35932.668
$ws.Range("L52").Value = 35932.668
$ws.Range("N52").Value = -36384.668
$ws.Range("H126").Value = 2815.8635
$ws.Range("I126").Value = 1156
$ws.Range("J126").Value = 7242.1665
$ws.Range("K126").Value = 3468
$ws.Range("L126").Value = 21726.4995
$ws.Range("M126").Value = -998
$ws.Range("N126").Value = -26666.4995
$ws.Range("H130").Value = 35799.6
$ws.Range("J130").Value = 35799.6
$ws.Range("L130").Value = 35799.6
$ws.Range("N130").Value = -45839.6
